$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 500.12903
$ws.Range("I33").Value = 185.84616
$ws.Range("J33").Value = 2134.4
$ws.Range("K33").Value = 185.84616
$ws.Range("L33").Value = 2134.4
$ws.Range("M33").Value = 43.15384
$ws.Range("N33").Value = -2592.4

$ws.Range("H76").Value = 66670370
$ws.Range("I76").Value = 90911816
$ws.Range("J76").Value = 6377
$ws.Range("K76").Value = 90911816
$ws.Range("L76").Value = 6377
$ws.Range("M76").Value = -90911501
$ws.Range("N76").Value = -7007

$ws.Range("H79").Value = 66670370
$ws.Range("I79").Value = 90911816
$ws.Range("J79").Value = 6377
$ws.Range("K79").Value = 90911816
$ws.Range("L79").Value = 6377
$ws.Range("M79").Value = -90910724
$ws.Range("N79").Value = -8561

$ws.Range("H123").Value = 34280
$ws.Range("J123").Value = 34280
$ws.Range("L123").Value = 34280
$ws.Range("N123").Value = -44080

$ws.Range("H125").Value = 58825030
$ws.Range("I125").Value = 83334520
$ws.Range("J125").Value = 2259.8
$ws.Range("K125").Value = 750010680
$ws.Range("L125").Value = 20338.2
$ws.Range("M125").Value = -750008220
$ws.Range("N125").Value = -25258.2

$ws.Range("H138").Value = 5156665.5
$ws.Range("I138").Value = 1192.1562
$ws.Range("J138").Value = 7694745
$ws.Range("K138").Value = 3576.4686
$ws.Range("L138").Value = 23084235
$ws.Range("M138").Value = 1563.5314
$ws.Range("N138").Value = -23094515

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2764.3635
$ws.Range("I63").Value = 2764.3635
$ws.Range("K63").Value = 2764.3635
$ws.Range("M63").Value = -2078.3635

$ws.Range("H66").Value = 2764.3635
$ws.Range("I66").Value = 2764.3635
$ws.Range("K66").Value = 13821.8175
$ws.Range("M66").Value = -10389.8175

$ws.Range("H108").Value = 19684
$ws.Range("J108").Value = 19684
$ws.Range("L108").Value = 19684
$ws.Range("N108").Value = -27364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 443160.66
$ws.Range("I105").Value = 531552.8
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 531552.8
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = -529805.8
$ws.Range("N105").Value = -4694

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24047.652
$ws.Range("I31").Value = 25443.707
$ws.Range("K31").Value = 25443.707
$ws.Range("M31").Value = -25148.707

$ws.Range("H34").Value = 24047.652
$ws.Range("I34").Value = 25443.707
$ws.Range("K34").Value = 25443.707
$ws.Range("M34").Value = -25241.707

$ws.Range("H36").Value = 75053
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 75053
$ws.Range("K36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("M36").Value = 75053
$ws.Range("N36").Value = -75829

$ws.Range("H40").Value = 75053
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 75053
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = 75053
$ws.Range("N40").Value = -75373

$ws.Range("H62").Value = 3004.1667
$ws.Range("I62").Value = 2959.0908
$ws.Range("K62").Value = 2959.0908
$ws.Range("M62").Value = -2335.0908

$ws.Range("H65").Value = 3004.1667
$ws.Range("I65").Value = 2959.0908
$ws.Range("K65").Value = 14795.454
$ws.Range("M65").Value = -11675.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2551
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 5002
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 30012
$ws.Range("M2").Value = -487
$ws.Range("N2").Value = -30238

$ws.Range("H34").Value = 12820611
$ws.Range("J34").Value = 12820611
$ws.Range("L34").Value = 38461833
$ws.Range("N34").Value = -38462001

$ws.Range("H39").Value = 10833.917
$ws.Range("J39").Value = 10833.917
$ws.Range("L39").Value = 32501.751
$ws.Range("N39").Value = -33089.751

$ws.Range("H55").Value = 91669656
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 91669656
$ws.Range("K55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").Value = 275008968
$ws.Range("N55").Value = -275009322

$ws.Range("H105").Value = 12866.333
$ws.Range("J105").Value = 12866.333
$ws.Range("L105").Value = 38598.999
$ws.Range("N105").Value = -43840.999

$ws.Range("H131").Value = 17361916
$ws.Range("I131").Value = 439.0909
$ws.Range("J131").Value = 20492674
$ws.Range("K131").Value = 1317.2727
$ws.Range("L131").Value = 61478022
$ws.Range("M131").Value = 3722.7273
$ws.Range("N131").Value = -61488102

$ws.Range("H134").Value = 4235.225
$ws.Range("I134").Value = 1284.238
$ws.Range("J134").Value = 7496.8423
$ws.Range("K134").Value = 3852.714
$ws.Range("L134").Value = 22490.5269
$ws.Range("M134").Value = 1217.286
$ws.Range("N134").Value = -32630.5269

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26372

$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 25000
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81864

$ws.Range("H70").Value = 5887094
$ws.Range("I70").Value = 11115322
$ws.Range("J70").Value = 5337
$ws.Range("K70").Value = 11115322
$ws.Range("L70").Value = 5337
$ws.Range("M70").Value = -11115052
$ws.Range("N70").Value = -5877

$ws.Range("H73").Value = 5887094
$ws.Range("I73").Value = 11115322
$ws.Range("J73").Value = 5337
$ws.Range("K73").Value = 11115322
$ws.Range("L73").Value = 5337
$ws.Range("M73").Value = -11114386
$ws.Range("N73").Value = -7209

$ws.Range("H80").Value = 6029
$ws.Range("I80").Value = 2615
$ws.Range("J80").Value = 9443
$ws.Range("K80").Value = 2615
$ws.Range("L80").Value = 9443
$ws.Range("M80").Value = -1617
$ws.Range("N80").Value = -11439

$ws.Range("H83").Value = 6029
$ws.Range("I83").Value = 2615
$ws.Range("J83").Value = 9443
$ws.Range("K83").Value = 13075
$ws.Range("L83").Value = 47215
$ws.Range("M83").Value = -8083
$ws.Range("N83").Value = -57199

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 100000
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101498

$ws.Range("H66").Value = 100000
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -307488
